$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.034649
$ws.Range("H2").Value = 3.103947
$ws.Range("I2").Value = 0.02307585235784855
$ws.Range("J2").Value = 0.02307585235784855
$ws.Range("M2").Value = 179.7005413333333
$ws.Range("N2").Value = 539.101624
$ws.Range("O2").Value = 0.7012656334041908
$ws.Range("P2").Value = 0.7012656334041907
$ws.Range("Q2").Value = 185.926985389992
$ws.Range("R2").Value = 1673.342868509928
$ws.Range("S2").Value = 0.01618230222006825
$ws.Range("T2").Value = 0.01618230222006825
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.034649
$ws.Range("H3").Value = 3.103947
$ws.Range("I3").Value = 0.02307585235784855
$ws.Range("J3").Value = 0.02307585235784855
$ws.Range("O3").Value = 0.05908927597267952
$ws.Range("P3").Value = 0.05908927597267952
$ws.Range("Q3").Value = 15.66637580277
$ws.Range("R3").Value = 140.99738222493
$ws.Range("S3").Value = 0.001363535408277721
$ws.Range("T3").Value = 0.001363535408277721
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.034649
$ws.Range("H4").Value = 3.103947
$ws.Range("I4").Value = 0.02307585235784855
$ws.Range("J4").Value = 0.02307585235784855
$ws.Range("M4").Value = 36.14947766666667
$ws.Range("N4").Value = 108.448433
$ws.Range("O4").Value = 0.1410701724382803
$ws.Range("P4").Value = 0.1410701724382803
$ws.Range("Q4").Value = 37.402020918339
$ws.Range("R4").Value = 336.618188265051
$ws.Range("S4").Value = 0.003255314471281992
$ws.Range("T4").Value = 0.003255314471281992
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.034649
$ws.Range("H5").Value = 3.103947
$ws.Range("I5").Value = 0.02307585235784855
$ws.Range("J5").Value = 0.02307585235784855
$ws.Range("M5").Value = 25.25999466666667
$ws.Range("N5").Value = 75.779984
$ws.Range("O5").Value = 0.09857491818484938
$ws.Range("P5").Value = 0.09857491818484938
$ws.Range("Q5").Value = 26.135228221872
$ws.Range("R5").Value = 235.217053996848
$ws.Range("S5").Value = 0.002274700258220585
$ws.Range("T5").Value = 0.002274700258220585
$ws.Range("I6").Value = 0.9141262989281937
$ws.Range("J6").Value = 0.9141262989281937
$ws.Range("M6").Value = 179.7005413333333
$ws.Range("N6").Value = 539.101624
$ws.Range("O6").Value = 0.7012656334041908
$ws.Range("P6").Value = 0.7012656334041907
$ws.Range("Q6").Value = 7365.30743869241
$ws.Range("R6").Value = 66287.76694823169
$ws.Range("S6").Value = 0.6410453580293084
$ws.Range("T6").Value = 0.6410453580293083
$ws.Range("I7").Value = 0.9141262989281937
$ws.Range("J7").Value = 0.9141262989281937
$ws.Range("O7").Value = 0.05908927597267952
$ws.Range("P7").Value = 0.05908927597267952
$ws.Range("S7").Value = 0.05401506115125218
$ws.Range("T7").Value = 0.05401506115125218
$ws.Range("I8").Value = 0.9141262989281937
$ws.Range("J8").Value = 0.9141262989281937
$ws.Range("M8").Value = 36.14947766666667
$ws.Range("N8").Value = 108.448433
$ws.Range("O8").Value = 0.1410701724382803
$ws.Range("P8").Value = 0.1410701724382803
$ws.Range("Q8").Value = 1481.642819702275
$ws.Range("R8").Value = 13334.78537732047
$ws.Range("S8").Value = 0.1289559546201673
$ws.Range("T8").Value = 0.1289559546201673
$ws.Range("I9").Value = 0.9141262989281937
$ws.Range("J9").Value = 0.9141262989281937
$ws.Range("M9").Value = 25.25999466666667
$ws.Range("N9").Value = 75.779984
$ws.Range("O9").Value = 0.09857491818484938
$ws.Range("P9").Value = 0.09857491818484938
$ws.Range("Q9").Value = 1035.320345945001
$ws.Range("R9").Value = 9317.883113505008
$ws.Range("S9").Value = 0.09010992512746586
$ws.Range("T9").Value = 0.09010992512746586
$ws.Range("G10").Value = 2.775347666666667
$ws.Range("H10").Value = 8.326043
$ws.Range("I10").Value = 0.06189878209682655
$ws.Range("J10").Value = 0.06189878209682655
$ws.Range("M10").Value = 179.7005413333333
$ws.Range("N10").Value = 539.101624
$ws.Range("O10").Value = 0.7012656334041908
$ws.Range("P10").Value = 0.7012656334041907
$ws.Range("Q10").Value = 498.7314780882036
$ws.Range("R10").Value = 4488.583302793832
$ws.Range("S10").Value = 0.04340748863407905
$ws.Range("T10").Value = 0.04340748863407905
$ws.Range("G11").Value = 2.775347666666667
$ws.Range("H11").Value = 8.326043
$ws.Range("I11").Value = 0.06189878209682655
$ws.Range("J11").Value = 0.06189878209682655
$ws.Range("O11").Value = 0.05908927597267952
$ws.Range("P11").Value = 0.05908927597267952
$ws.Range("Q11").Value = 42.02356502479667
$ws.Range("R11").Value = 378.21208522317
$ws.Range("S11").Value = 0.003657554217692138
$ws.Range("T11").Value = 0.003657554217692138
$ws.Range("G12").Value = 2.775347666666667
$ws.Range("H12").Value = 8.326043
$ws.Range("I12").Value = 0.06189878209682655
$ws.Range("J12").Value = 0.06189878209682655
$ws.Range("M12").Value = 36.14947766666667
$ws.Range("N12").Value = 108.448433
$ws.Range("O12").Value = 0.1410701724382803
$ws.Range("P12").Value = 0.1410701724382803
$ws.Range("Q12").Value = 100.3273684934021
$ws.Range("R12").Value = 902.9463164406191
$ws.Range("S12").Value = 0.008732071864118857
$ws.Range("T12").Value = 0.008732071864118857
$ws.Range("G13").Value = 2.775347666666667
$ws.Range("H13").Value = 8.326043
$ws.Range("I13").Value = 0.06189878209682655
$ws.Range("J13").Value = 0.06189878209682655
$ws.Range("M13").Value = 25.25999466666667
$ws.Range("N13").Value = 75.779984
$ws.Range("O13").Value = 0.09857491818484938
$ws.Range("P13").Value = 0.09857491818484938
$ws.Range("Q13").Value = 70.10526725814579
$ws.Range("R13").Value = 630.9474053233121
$ws.Range("S13").Value = 0.006101667380936497
$ws.Range("T13").Value = 0.006101667380936497
$ws.Range("G14").Value = 0.04031133333333333
$ws.Range("H14").Value = 0.120934
$ws.Range("I14").Value = 0.0008990666171310454
$ws.Range("J14").Value = 0.0008990666171310454
$ws.Range("M14").Value = 179.7005413333333
$ws.Range("N14").Value = 539.101624
$ws.Range("O14").Value = 0.7012656334041908
$ws.Range("P14").Value = 0.7012656334041907
$ws.Range("Q14").Value = 7.243968421868445
$ws.Range("R14").Value = 65.195715796816
$ws.Range("S14").Value = 0.0006304845207349656
$ws.Range("T14").Value = 0.0006304845207349655
$ws.Range("G15").Value = 0.04031133333333333
$ws.Range("H15").Value = 0.120934
$ws.Range("I15").Value = 0.0008990666171310454
$ws.Range("J15").Value = 0.0008990666171310454
$ws.Range("O15").Value = 0.05908927597267952
$ws.Range("P15").Value = 0.05908927597267952
$ws.Range("Q15").Value = 0.6103833252733333
$ws.Range("R15").Value = 5.49344992746
$ws.Range("S15").Value = 0.00005312519545747974
$ws.Range("T15").Value = 0.00005312519545747974
$ws.Range("G16").Value = 0.04031133333333333
$ws.Range("H16").Value = 0.120934
$ws.Range("I16").Value = 0.0008990666171310454
$ws.Range("J16").Value = 0.0008990666171310454
$ws.Range("M16").Value = 36.14947766666667
$ws.Range("N16").Value = 108.448433
$ws.Range("O16").Value = 0.1410701724382803
$ws.Range("P16").Value = 0.1410701724382803
$ws.Range("Q16").Value = 1.457233644046889
$ws.Range("R16").Value = 13.115102796422
$ws.Range("S16").Value = 0.0001268314827121779
$ws.Range("T16").Value = 0.0001268314827121779
$ws.Range("G17").Value = 0.04031133333333333
$ws.Range("H17").Value = 0.120934
$ws.Range("I17").Value = 0.0008990666171310454
$ws.Range("J17").Value = 0.0008990666171310454
$ws.Range("M17").Value = 25.25999466666667
$ws.Range("N17").Value = 75.779984
$ws.Range("O17").Value = 0.09857491818484938
$ws.Range("P17").Value = 0.09857491818484938
$ws.Range("Q17").Value = 1.018264065006222
$ws.Range("R17").Value = 9.164376585055999
$ws.Range("S17").Value = 0.00008862541822642211
$ws.Range("T17").Value = 0.00008862541822642211
